$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-458) holds the "Förändrad" (last changed) date, stored
# as serial date 45180 (2023-09-11). Bump it by one day to 45181
# (2023-09-12) for every data row, matching the diff.
$ws.Range("C2:C458").Value = 45181
